# Atualizando o arquivo XLSX
# Apply updated odds values to row 2 (match: San Martin T. vs San Telmo)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = 6        # Odd_A_FT: 5.75 -> 6
$ws.Range("J2").Value = 2.5      # Odd_H_HT: 2.6 -> 2.5
$ws.Range("N2").Value = 5        # Odd_Under05_FT: 4.75 -> 5
$ws.Range("Q2").Value = 3.5      # Odd_Over25_FT: 3.4 -> 3.5
$ws.Range("R2").Value = 1.3      # Odd_Under25_FT: 1.33 -> 1.3
$ws.Range("X2").Value = 6        # Odd_CS_2-0: 6.5 -> 6
$ws.Range("AN2").Value = 3.25    # Odd_CS_1-0_HT: 3.4 -> 3.25
$ws.Range("AV2").Value = 7.5     # Odd_CS_0-1_HT: 7 -> 7.5
